$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1502.0834
$ws.Range("J70").Value = 1813.4
$ws.Range("L70").Value = 5440.200000000001
$ws.Range("N70").Value = -5980.200000000001
$ws.Range("H73").Value = 1502.0834
$ws.Range("J73").Value = 1813.4
$ws.Range("L73").Value = 5440.200000000001
$ws.Range("N73").Value = -7312.200000000001
$ws.Range("J86").Value = 3919
$ws.Range("L86").Value = 3919
$ws.Range("N86").Value = -6165
$ws.Range("J89").Value = 3919
$ws.Range("L89").Value = 19595
$ws.Range("N89").Value = -30827
$ws.Range("H138").Value = 5232.1035
$ws.Range("I138").Value = 10594.267
$ws.Range("K138").Value = 31782.801
$ws.Range("M138").Value = -26642.801

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1944.5
$ws.Range("I2").Value = 1333.4
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 1333.4
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -1220.4
$ws.Range("N2").Value = -5226
$ws.Range("H45").Value = 504443.5
$ws.Range("I45").Value = 504443.5
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 504443.5
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -504066.5
$ws.Range("N45").ClearContents()
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H116").Value = 1944.5
$ws.Range("I116").Value = 1333.4
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 1333.4
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = 960.5999999999999
$ws.Range("N116").Value = -9588
$ws.Range("H132").Value = 3814.84
$ws.Range("J132").Value = 4010.3333
$ws.Range("L132").Value = 12030.9999
$ws.Range("N132").Value = -17090.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1944.5
$ws.Range("I3").Value = 1333.4
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 1333.4
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -1219.4
$ws.Range("N3").Value = -5228
$ws.Range("H134").Value = 25716696
$ws.Range("I134").Value = 2233.5652
$ws.Range("K134").Value = 6700.6956
$ws.Range("M134").Value = -4165.6956

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2436.2
$ws.Range("I31").Value = 2035.0769
$ws.Range("K31").Value = 2035.0769
$ws.Range("M31").Value = -1740.0769
$ws.Range("H34").Value = 2436.2
$ws.Range("I34").Value = 2035.0769
$ws.Range("K34").Value = 2035.0769
$ws.Range("M34").Value = -1833.0769
$ws.Range("H58").Value = 2416.6897
$ws.Range("I58").Value = 2510.3
$ws.Range("K58").Value = 2510.3
$ws.Range("M58").Value = -2307.3
$ws.Range("H132").Value = 29952.139
$ws.Range("I132").Value = 43688.582
$ws.Range("K132").Value = 131065.746
$ws.Range("M132").Value = -128535.746
$ws.Range("H134").Value = 3069.9285
$ws.Range("I134").Value = 2688.75
$ws.Range("J134").Value = 3578.1667
$ws.Range("K134").Value = 8066.25
$ws.Range("L134").Value = 10734.5001
$ws.Range("M134").Value = -5531.25
$ws.Range("N134").Value = -15804.5001
$ws.Range("H136").Value = 2416.6897
$ws.Range("I136").Value = 2510.3
$ws.Range("K136").Value = 7530.900000000001
$ws.Range("M136").Value = -4980.900000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 17105.875
$ws.Range("J62").Value = 19641.334
$ws.Range("L62").Value = 58924.00199999999
$ws.Range("N62").Value = -60296.00199999999
$ws.Range("H65").Value = 17105.875
$ws.Range("J65").Value = 19641.334
$ws.Range("L65").Value = 176772.006
$ws.Range("N65").Value = -183636.006
$ws.Range("H120").Value = 22548.8
$ws.Range("J120").Value = 24498.666
$ws.Range("L120").Value = 73495.99800000001
$ws.Range("N120").Value = -83171.99800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 65160.434
$ws.Range("I80").Value = 11477.929
$ws.Range("K80").Value = 11477.929
$ws.Range("M80").Value = -10479.929
$ws.Range("H83").Value = 65160.434
$ws.Range("I83").Value = 11477.929
$ws.Range("K83").Value = 57389.645
$ws.Range("M83").Value = -52397.645
$ws.Range("H102").Value = 15626909
$ws.Range("I102").Value = 19232458
$ws.Range("K102").Value = 19232458
$ws.Range("M102").Value = -19230836
$ws.Range("H113").Value = 1976.75
$ws.Range("I113").Value = 1752.625
$ws.Range("K113").Value = 1752.625
$ws.Range("M113").Value = 417.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 948.3333
$ws.Range("I16").Value = 889.3333
$ws.Range("K16").Value = 889.3333
$ws.Range("M16").Value = -719.3333
$ws.Range("H22").Value = 4898.524
$ws.Range("I22").Value = 2282.2
$ws.Range("K22").Value = 2282.2
$ws.Range("M22").Value = -1987.2
$ws.Range("H27").Value = 4898.524
$ws.Range("I27").Value = 2282.2
$ws.Range("K27").Value = 2282.2
$ws.Range("M27").Value = -2175.2
$ws.Range("H132").Value = 3032.8108
$ws.Range("I132").Value = 2778.25
$ws.Range("K132").Value = 8334.75
$ws.Range("M132").Value = -5804.75
$ws.Range("H136").Value = 1971.2699
$ws.Range("I136").Value = 3248.5
$ws.Range("J136").Value = 1670.7451
$ws.Range("K136").Value = 9745.5
$ws.Range("L136").Value = 5012.2353
$ws.Range("M136").Value = -7195.5
$ws.Range("N136").Value = -10112.2353

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 42399.5
$ws.Range("I70").Value = 34799
$ws.Range("K70").Value = 34799
$ws.Range("M70").Value = -34484
$ws.Range("H73").Value = 42399.5
$ws.Range("I73").Value = 34799
$ws.Range("K73").Value = 34799
$ws.Range("M73").Value = -33707
$ws.Range("H81").Value = 75457.78999999999
$ws.Range("J81").Value = 203159.4
$ws.Range("L81").Value = 406318.8
$ws.Range("N81").Value = -408440.8
$ws.Range("H84").Value = 75457.78999999999
$ws.Range("J84").Value = 203159.4
$ws.Range("L84").Value = 2031594
$ws.Range("N84").Value = -2042202
$ws.Range("H107").Value = 840844.7
$ws.Range("I107").Value = 540.03705
$ws.Range("J107").Value = 4082019.8
$ws.Range("K107").Value = 1620.11115
$ws.Range("L107").Value = 12246059.4
$ws.Range("M107").Value = 299.8888499999998
$ws.Range("N107").Value = -12249899.4
$ws.Range("H113").Value = 608.1111
$ws.Range("I113").Value = 495.5
$ws.Range("J113").Value = 1002.25
$ws.Range("K113").Value = 1486.5
$ws.Range("L113").Value = 3006.75
$ws.Range("M113").Value = 683.5
$ws.Range("N113").Value = -7346.75

